$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 150, shifting old rows 150:153 down to 151:154
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with data (mirrors the row that used to
# be directly above it, i.e. row 149 / row 151 after the shift)
$ws.Range("A150").Value = 1
$ws.Range("B150").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C150").Value = "Arica y Parinacota"
$ws.Range("D150").Value = 44939
$ws.Range("E150").Value = 15
$ws.Range("F150").Value = 100112042
$ws.Range("G150").Value = "Locoto"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 150
$ws.Range("K150").Value = 19000
$ws.Range("L150").Value = 20000
$ws.Range("M150").Value = 19400
$ws.Range("N150").Value = "$/caja 20 kilos"
$ws.Range("O150").Value = "Región de Arica y Parinacota"
$ws.Range("P150").Value = 970
$ws.Range("Q150").Value = 20
$ws.Range("R150").Value = "Hortaliza"
